$d = $word.ActiveDocument

# Remove every comment in the document (both comments were resolved /
# no longer needed per the updated guidelines). Deleting each Comment
# object removes its commentRangeStart/commentRangeEnd/commentReference
# markers from the body text as well as the comment entry itself.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments($i).Delete()
}
